$d = $word.ActiveDocument

# Locate the paragraph introducing the constellation ("Sternbildes Perseus ...")
$target = $null
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Mach mit an einer weltweiten Kampagne*") {
        $target = $p
        break
    }
}

$newText = "Mach mit an einer weltweiten Kampagne, die schwächsten sichtbaren Sterne zu beobachten und aufzuzeichnen, um die Lichtverschmutzung an einem Ort zu messen. Durch das Auffinden und Beobachten des Sternbild Herkules am Nachthimmel und den Vergleich mit den Helligkeitskarten, lernen Menschen auf der ganzen Erde, wie die Lichter in ihrer Gemeinde zur Lichtverschmutzung beitragen. Dein Beitrag zur Online-Datenbank beschreibt den sichtbaren Nachthimmel."

# Select the paragraph's text (excluding the trailing paragraph mark), remove the
# four differently-formatted runs that currently make it up, and retype the whole
# sentence fresh as a single, unformatted run.
$r = $target.Range
$r.End = $r.End - 1
$r.Delete()
$r.InsertAfter($newText)
